# Helper: assign a value to a cell while forcing it to stay plain text,
# matching the source workbook (all data cells are inline/shared strings,
# e.g. "207.08" must not silently become the number 207.08, and "0.0590"
# must not lose its trailing zero). We temporarily switch the cell to the
# Text number format for the assignment, then restore the original style
# so the cell formatting/style is left exactly as it was.
function Set-TextValue {
    param($ws, $row, $col, $val)
    $cell = $ws.Cells.Item($row, $col)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $savedStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws 2 4 "27.684.52"
Set-TextValue $ws 2 5 "  -0.73%  "

Set-TextValue $ws 3 4 "1.583.44"
Set-TextValue $ws 3 5 "  -3.05%  "

Set-TextValue $ws 4 5 "  +0.03%  "

Set-TextValue $ws 5 4 "207.08"
Set-TextValue $ws 5 5 "  -2.14%  "

Set-TextValue $ws 6 4 "0.503"
Set-TextValue $ws 6 5 "  -3.26%  "

Set-TextValue $ws 7 5 "  +0.05%  "

Set-TextValue $ws 8 4 "22.31"
Set-TextValue $ws 8 5 "  -4.51%  "

Set-TextValue $ws 9 5 "  -1.44%  "

Set-TextValue $ws 10 4 "0.0590"
Set-TextValue $ws 10 5 "  -3.57%  "

Set-TextValue $ws 11 4 "0.0868"
Set-TextValue $ws 11 5 "  -1.70%  "

Set-TextValue $ws 12 4 "1.806.76"
Set-TextValue $ws 12 5 "  -3.13%  "

Set-TextValue $ws 13 4 "1.575.63"
Set-TextValue $ws 13 5 "  -3.62%  "

Set-TextValue $ws 14 4 "3.86"
Set-TextValue $ws 14 5 "  -4.09%  "

Set-TextValue $ws 15 4 "0.530"
Set-TextValue $ws 15 5 "  -5.79%  "

Set-TextValue $ws 16 4 "27.634.61"
Set-TextValue $ws 16 5 "  -0.94%  "

Set-TextValue $ws 17 4 "62.86"
Set-TextValue $ws 17 5 "  -3.82%  "

Set-TextValue $ws 18 4 "217.69"
Set-TextValue $ws 18 5 "  -5.01%  "

Set-TextValue $ws 19 2 "ShibaInu"
Set-TextValue $ws 19 3 "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$u2083 = [char]0x2083
Set-TextValue $ws 19 4 "0.0${u2083}0693"
Set-TextValue $ws 19 5 "  -3.64%  "

Set-TextValue $ws 20 2 "Chainlink"
Set-TextValue $ws 20 3 "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws 20 4 "7.33"
Set-TextValue $ws 20 5 "  -4.38%  "

Set-TextValue $ws 21 5 "  +0.11%  "

Set-TextValue $ws 22 4 "4.15"
Set-TextValue $ws 22 5 "  -4.28%  "

Set-TextValue $ws 23 4 "9.53"
Set-TextValue $ws 23 5 "  -5.07%  "

Set-TextValue $ws 24 5 "  -4.37%  "

Set-TextValue $ws 25 4 "153.44"
Set-TextValue $ws 25 5 "  -1.38%  "

Set-TextValue $ws 26 5 "  +0.04%  "

Set-TextValue $ws 27 4 "6.69"
Set-TextValue $ws 27 5 "  -2.64%  "

Set-TextValue $ws 28 4 "15.07"
Set-TextValue $ws 28 5 "  -3.09%  "

Set-TextValue $ws 29 5 "  -4.36%  "

Set-TextValue $ws 30 4 "1.16"
Set-TextValue $ws 30 5 "  -1.94%  "

Set-TextValue $ws 31 5 "  -3.79%  "

Set-TextValue $ws 32 5 "  -5.22%  "

Set-TextValue $ws 33 4 "1.374.86"
Set-TextValue $ws 33 5 "  -1.39%  "

Set-TextValue $ws 34 4 "2.95"
Set-TextValue $ws 34 5 "  -5.15%  "

Set-TextValue $ws 35 4 "1.52"
Set-TextValue $ws 35 5 "  -4.95%  "

Set-TextValue $ws 36 4 "0.969"
Set-TextValue $ws 36 5 "  -5.19%  "

Set-TextValue $ws 37 5 "  -2.15%  "

Set-TextValue $ws 38 4 "0.0165"
Set-TextValue $ws 38 5 "  -3.38%  "

Set-TextValue $ws 39 4 "0.537"
Set-TextValue $ws 39 5 "  -4.04%  "

Set-TextValue $ws 40 4 "0.815"
Set-TextValue $ws 40 5 "  -3.90%  "

Set-TextValue $ws 41 5 "  +0.10%  "

Set-TextValue $ws 42 4 "0.978"
Set-TextValue $ws 42 5 "  -3.51%  "

Set-TextValue $ws 43 4 "1.78"
Set-TextValue $ws 43 5 "  -2.45%  "

Set-TextValue $ws 44 5 "  +1.62%  "

Set-TextValue $ws 45 4 "63.64"
Set-TextValue $ws 45 5 "  -3.45%  "

Set-TextValue $ws 46 5 "  -3.45%  "

Set-TextValue $ws 47 4 "1.718.92"
Set-TextValue $ws 47 5 "  -3.09%  "

Set-TextValue $ws 48 4 "87.50"
Set-TextValue $ws 48 5 "  -1.64%  "

$u2086 = [char]0x2086
Set-TextValue $ws 49 4 "0.0${u2086}0101"
Set-TextValue $ws 49 5 "  -2.28%  "

Set-TextValue $ws 50 4 "0.0974"
Set-TextValue $ws 50 5 "  -4.41%  "

Set-TextValue $ws 51 5 "  -1.53%  "

